$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: kamal / IT / 14-Jan-2019 / 70,000 LKR
$ws.Range("A2").Value = "kamal"
$ws.Range("B2").Value = "IT"
$ws.Range("C2").Value = 43479
$ws.Range("D2").Value = "70,000 LKR"

# Row 3: saman / flightops / 12-Dec-2025 / 20,000 LKR
$ws.Range("A3").Value = "saman"
$ws.Range("B3").Value = "flightops"
$ws.Range("C3").Value = 46003
$ws.Range("D3").Value = "20,000 LKR"
